$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsALC.Range("H28").Value = 1601.3334
$wsALC.Range("J28").Value = 4375.8
$wsALC.Range("L28").Value = 4375.8
$wsALC.Range("N28").Value = -5345.8

$wsALC.Range("H32").Value = 1095.4
$wsALC.Range("I32").Value = 1000
$wsALC.Range("J32").Value = 1119.25
$wsALC.Range("K32").Value = 1000
$wsALC.Range("L32").Value = 1119.25
$wsALC.Range("M32").Value = -674
$wsALC.Range("N32").Value = -1771.25

$wsALC.Range("H33").Value = 123.125
$wsALC.Range("I33").Value = 105
$wsALC.Range("K33").Value = 105
$wsALC.Range("M33").Value = 124

$wsALC.Range("H62").Value = 4146.2173
$wsALC.Range("I62").Value = 3719.3684
$wsALC.Range("K62").Value = 3719.3684
$wsALC.Range("M62").Value = -3095.3684

$wsALC.Range("H65").Value = 4146.2173
$wsALC.Range("I65").Value = 3719.3684
$wsALC.Range("K65").Value = 18596.842
$wsALC.Range("M65").Value = -15476.842

$wsALC.Range("H98").Value = 1490.6957
$wsALC.Range("I98").Value = 877.1177
$wsALC.Range("J98").Value = 3229.1667
$wsALC.Range("K98").Value = 877.1177
$wsALC.Range("L98").Value = 3229.1667
$wsALC.Range("M98").Value = 620.8823
$wsALC.Range("N98").Value = -6225.1667

$wsALC.Range("H122").Value = 1490.6957
$wsALC.Range("I122").Value = 877.1177
$wsALC.Range("J122").Value = 3229.1667
$wsALC.Range("K122").Value = 2631.3531
$wsALC.Range("L122").Value = 9687.500100000001
$wsALC.Range("M122").Value = -181.3531000000003
$wsALC.Range("N122").Value = -14587.5001

$wsALC.Range("H132").Value = 1830.9762
$wsALC.Range("I132").Value = 1830.9762
$wsALC.Range("K132").Value = 5492.9286
$wsALC.Range("M132").Value = -2962.9286

$wsALC.Range("H137").Value = 1916.72
$wsALC.Range("I137").Value = 1707.421
$wsALC.Range("J137").Value = 2579.5
$wsALC.Range("K137").Value = 5122.263
$wsALC.Range("L137").Value = 7738.5
$wsALC.Range("M137").Value = -2572.263
$wsALC.Range("N137").Value = -12838.5

$wsALC.Range("H138").Value = 2360.422
$wsALC.Range("J138").Value = 2092.9355
$wsALC.Range("L138").Value = 6278.806500000001
$wsALC.Range("N138").Value = -16558.8065

$wsARM = $wb.Worksheets.Item("ARM")
$wsARM.Range("H45").Value = 6580.731
$wsARM.Range("I45").Value = 9798.77
$wsARM.Range("K45").Value = 9798.77
$wsARM.Range("M45").Value = -9421.77

$wsARM.Range("H74").Value = 1561.5741
$wsARM.Range("I74").Value = 1398.8
$wsARM.Range("K74").Value = 1398.8
$wsARM.Range("M74").Value = -524.8

$wsARM.Range("H77").Value = 1561.5741
$wsARM.Range("I77").Value = 1398.8
$wsARM.Range("K77").Value = 6994
$wsARM.Range("M77").Value = -2626

$wsARM.Range("H132").Value = 3090.9688
$wsARM.Range("I132").Value = 1391.6086
$wsARM.Range("K132").Value = 4174.825800000001
$wsARM.Range("M132").Value = -1644.825800000001

$wsBSM = $wb.Worksheets.Item("BSM")
$wsBSM.Range("H20").Value = 13851.765
$wsBSM.Range("I20").Value = 12311.75
$wsBSM.Range("J20").Value = 17547.8
$wsBSM.Range("K20").Value = 12311.75
$wsBSM.Range("L20").Value = 17547.8
$wsBSM.Range("M20").Value = -12064.75
$wsBSM.Range("N20").Value = -18041.8

$wsBSM.Range("H94").Value = 4101.913
$wsBSM.Range("I94").Value = 3944.4211
$wsBSM.Range("J94").Value = 4850
$wsBSM.Range("K94").Value = 3944.4211
$wsBSM.Range("L94").Value = 4850
$wsBSM.Range("M94").Value = -3493.4211
$wsBSM.Range("N94").Value = -5752

$wsBSM.Range("H131").Value = 0
$wsBSM.Range("J131").Value = 0
$wsBSM.Range("L131").Value = 0
$wsBSM.Range("N131").Value = $null

$wsCRP = $wb.Worksheets.Item("CRP")
$wsCRP.Range("H132").Value = 2483.392
$wsCRP.Range("I132").Value = 2499.0571
$wsCRP.Range("K132").Value = 7497.1713
$wsCRP.Range("M132").Value = -4967.1713

$wsCUL = $wb.Worksheets.Item("CUL")
$wsCUL.Range("H24").Value = 116.666664
$wsCUL.Range("I24").Value = 116.666664
$wsCUL.Range("K24").Value = 349.999992
$wsCUL.Range("M24").Value = -119.999992

$wsGSM = $wb.Worksheets.Item("GSM")
$wsGSM.Range("H70").Value = 5427.636
$wsGSM.Range("I70").Value = 5563
$wsGSM.Range("K70").Value = 5563
$wsGSM.Range("M70").Value = -5293

$wsGSM.Range("H73").Value = 5427.636
$wsGSM.Range("I73").Value = 5563
$wsGSM.Range("K73").Value = 5563
$wsGSM.Range("M73").Value = -4627

$wsGSM.Range("H97").Value = 33535.047
$wsGSM.Range("I97").Value = 50374.145
$wsGSM.Range("J97").Value = 4066.625
$wsGSM.Range("K97").Value = 50374.145
$wsGSM.Range("L97").Value = 4066.625
$wsGSM.Range("M97").Value = -49878.145
$wsGSM.Range("N97").Value = -5058.625

$wsGSM.Range("H129").Value = 50000
$wsGSM.Range("I129").Value = 50000
$wsGSM.Range("J129").Value = 0
$wsGSM.Range("K129").Value = 50000
$wsGSM.Range("L129").Value = 0
$wsGSM.Range("M129").Value = -45000
$wsGSM.Range("N129").Value = $null

$wsLTW = $wb.Worksheets.Item("LTW")
$wsLTW.Range("H7").Value = 5975.891
$wsLTW.Range("J7").Value = 5242.1924
$wsLTW.Range("L7").Value = 5242.1924
$wsLTW.Range("N7").Value = -5466.1924

$wsLTW.Range("H22").Value = 2989.8
$wsLTW.Range("I22").Value = 2278.5833
$wsLTW.Range("J22").Value = 4056.625
$wsLTW.Range("K22").Value = 2278.5833
$wsLTW.Range("L22").Value = 4056.625
$wsLTW.Range("M22").Value = -1983.5833
$wsLTW.Range("N22").Value = -4646.625

$wsLTW.Range("H27").Value = 2989.8
$wsLTW.Range("I27").Value = 2278.5833
$wsLTW.Range("J27").Value = 4056.625
$wsLTW.Range("K27").Value = 2278.5833
$wsLTW.Range("L27").Value = 4056.625
$wsLTW.Range("M27").Value = -2171.5833
$wsLTW.Range("N27").Value = -4270.625

$wsLTW.Range("H46").Value = 25547.777
$wsLTW.Range("I46").Value = 47946.555
$wsLTW.Range("J46").Value = 3149
$wsLTW.Range("K46").Value = 47946.555
$wsLTW.Range("L46").Value = 3149
$wsLTW.Range("M46").Value = -47758.555
$wsLTW.Range("N46").Value = -3525

$wsLTW.Range("H55").Value = 216
$wsLTW.Range("I55").Value = 254.15384
$wsLTW.Range("J55").Value = 154
$wsLTW.Range("K55").Value = 254.15384
$wsLTW.Range("L55").Value = 154
$wsLTW.Range("M55").Value = -81.15384
$wsLTW.Range("N55").Value = -500

$wsLTW.Range("H59").Value = 100000
$wsLTW.Range("J59").Value = 100000
$wsLTW.Range("L59").Value = 100000
$wsLTW.Range("N59").Value = -101308

$wsLTW.Range("H100").Value = 16292.37
$wsLTW.Range("I100").Value = 1726.7273
$wsLTW.Range("K100").Value = 1726.7273
$wsLTW.Range("M100").Value = -1185.7273

$wsLTW.Range("H126").Value = 5975.891
$wsLTW.Range("J126").Value = 5242.1924
$wsLTW.Range("L126").Value = 15726.5772
$wsLTW.Range("N126").Value = -20666.5772
